# refactor of VfM analysis: add "Project/Programme" header to the two
# existing sheets and add a new "Count" summary sheet (PVC total per
# category + category count) after Q4_19_20.

$wb = $excel.ActiveWorkbook

# --- Existing sheets: label column B header row 2 -------------------------
$ws1 = $wb.Worksheets.Item("Q1_20_21")
$ws1.Range("B2").Value = "Project/Programme"

$ws2 = $wb.Worksheets.Item("Q4_19_20")
$ws2.Range("B2").Value = "Project/Programme"

# --- New "Count" sheet, placed after Q4_19_20 ------------------------------
$wsCount = $wb.Worksheets.Add($null, $ws2)
$wsCount.Name = "Count"

# PVC total per category
$wsCount.Range("B2").Value = "PVC total per category"

$wsCount.Range("B3").Value = "Category"
$wsCount.Range("C3").Value = "Q1 20/21"
$wsCount.Range("D3").Value = "Q4 19/20"

$wsCount.Range("B4").Value = "Poor"
$wsCount.Range("C4").Value = 1172
$wsCount.Range("D4").Value = 1172

$wsCount.Range("B5").Value = "Low"
$wsCount.Range("C5").Value = 0
$wsCount.Range("D5").Value = 0

$wsCount.Range("B6").Value = "Medium"
$wsCount.Range("C6").Value = 2956
$wsCount.Range("D6").Value = 2831

$wsCount.Range("B7").Value = "High"
$wsCount.Range("C7").Value = 1761
$wsCount.Range("D7").Value = 1761

$wsCount.Range("B8").Value = "Very High"
$wsCount.Range("C8").Value = 2089
$wsCount.Range("D8").Value = 1481.6

$wsCount.Range("B9").Value = "Very High and Financially Positive"
$wsCount.Range("C9").Value = 0
$wsCount.Range("D9").Value = 0

$wsCount.Range("B10").Value = "Economically Positive"
$wsCount.Range("C10").Value = 0
$wsCount.Range("D10").Value = 0

$wsCount.Range("B11").Value = "Total"
$wsCount.Range("C11").Value = 7978
$wsCount.Range("D11").Value = 7245.6

# Category count
$wsCount.Range("B14").Value = "Category count"

$wsCount.Range("B15").Value = "Category"
$wsCount.Range("C15").Value = "Q1 20/21"
$wsCount.Range("D15").Value = "Q4 19/20"

$wsCount.Range("B16").Value = "Poor"
$wsCount.Range("C16").Value = 1
$wsCount.Range("D16").Value = 1

$wsCount.Range("B17").Value = "Low"
$wsCount.Range("C17").Value = 0
$wsCount.Range("D17").Value = 0

$wsCount.Range("B18").Value = "Medium"
$wsCount.Range("C18").Value = 1
$wsCount.Range("D18").Value = 1

$wsCount.Range("B19").Value = "High"
$wsCount.Range("C19").Value = 2
$wsCount.Range("D19").Value = 2

$wsCount.Range("B20").Value = "Very High"
$wsCount.Range("C20").Value = 1
$wsCount.Range("D20").Value = 1

$wsCount.Range("B21").Value = "Very High and Financially Positive"
$wsCount.Range("C21").Value = 0
$wsCount.Range("D21").Value = 0

$wsCount.Range("B22").Value = "Economically Positive"
$wsCount.Range("C22").Value = 0
$wsCount.Range("D22").Value = 0

$wsCount.Range("B23").Value = "Total"
$wsCount.Range("C23").Value = 5
$wsCount.Range("D23").Value = 5
